$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50 (shifts the existing rows 50-52 down to 51-53,
# carrying over their values and formatting, e.g. the date style on column D)
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly price record
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44706
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = "Poroto verde"
$ws.Range("H50").Value = "Magnum"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 22000
$ws.Range("L50").Value = 24000
$ws.Range("M50").Value = 23000
$ws.Range("N50").Value = "`$/saco 25 kilos"
$ws.Range("O50").Value = "Perú"
$ws.Range("P50").Value = 920
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
